# login_helper for code reusability.
# Appends a new batch of payment records (rows 102-128) to the payments sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ A=1; B="Sophia Brown"; C=32; D="2025-10-05 18:43:37" },
    @{ A=2; B="Sophia Miller"; C=4555; D="2025-10-05 18:43:40" },
    @{ A=3; B="John Miller"; C=5560; D="2025-10-05 18:43:43" },
    @{ A=4; B="Ava Miller"; C=3241; D="2025-10-05 18:43:46" },
    @{ A=5; B="Olivia Johnson"; C=3400; D="2025-10-05 18:43:49" },
    @{ A=6; B="James Johnson"; C=3400; D="2025-10-05 18:43:52" },
    @{ A=7; B="John Smith"; C=34; D="2025-10-05 18:43:54" },
    @{ A=8; B="Michael Garcia"; C=6774; D="2025-10-05 18:43:57" },
    @{ A=9; B="Emma Brown"; C=3400; D="2025-10-05 18:44:00" },
    @{ A=10; B="Ava Garcia"; C=32; D="2025-10-05 18:44:03" },
    @{ A=11; B="Ava Davis"; C=34234; D="2025-10-05 18:44:06" },
    @{ A=12; B="Ava Smith"; C=3100; D="2025-10-05 18:44:09" },
    @{ A=13; B="Olivia Brown"; C=1212; D="2025-10-05 18:44:12" },
    @{ A=14; B="Olivia Miller"; C=3241; D="2025-10-05 18:44:15" },
    @{ A=15; B="Ava Garcia"; C=3245; D="2025-10-05 18:44:18" },
    @{ A=16; B="Olivia Jones"; C=1212; D="2025-10-05 18:44:21" },
    @{ A=17; B="Michael Miller"; C=6774; D="2025-10-05 18:44:24" },
    @{ A=18; B="James Johnson"; C=6774; D="2025-10-05 18:44:27" },
    @{ A=19; B="Ava Smith"; C=4555; D="2025-10-05 18:44:30" },
    @{ A=20; B="James Garcia"; C=8987; D="2025-10-05 18:44:33" },
    @{ A=21; B="Ava Williams"; C=8987; D="2025-10-05 18:44:36" },
    @{ A=22; B="Ava Johnson"; C=34234; D="2025-10-05 18:44:39" },
    @{ A=23; B="John Johnson"; C=8987; D="2025-10-05 18:44:42" },
    @{ A=24; B="Sophia Brown"; C=3245; D="2025-10-05 18:44:45" },
    @{ A=25; B="Sophia Davis"; C=3245; D="2025-10-05 18:44:48" },
    @{ A=26; B="Michael Williams"; C=3400; D="2025-10-05 18:44:51" },
    @{ A=27; B="Michael Davis"; C=34; D="2025-10-05 18:44:54" }
)

$startRow = 102
foreach ($row in $newRows) {
    $ws.Cells.Item($startRow, 1).Value = $row.A
    $ws.Cells.Item($startRow, 2).Value = $row.B
    $ws.Cells.Item($startRow, 3).Value = $row.C
    $ws.Cells.Item($startRow, 4).Value = $row.D
    $startRow++
}
